$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-9, columns D, L, M, N, O, P, Q, R, S, T
# (Derived from a cyclic re-shuffle of the original row data.)

$data = @{
    2  = @{ D = 44742; L = "Segunda"; M = 100; N = 14000; O = 15000; P = 14500; Q = "`$/caja 18 kilos granel";     R = "Región de O'Higgins";  S = 806;   T = 18 }
    3  = @{ D = 45077; L = "Primera"; M = 140; N = 12000; O = 14000; P = 12857; Q = "`$/caja 12 kilos granel";     R = "Región de O'Higgins";  S = 12857; T = 1  }
    4  = @{ D = 45077; L = "Segunda"; M = 80;  N = 11000; O = 11000; P = 11000; Q = "`$/caja 12 kilos granel";     R = "Región de O'Higgins";  S = 11000; T = 1  }
    5  = @{ D = 44707; L = "Primera"; M = 60;  N = 12000; O = 13000; P = 12500; Q = "`$/caja 12 kilos empedrada";  R = "Provincia de Curicó";  S = 1042;  T = 12 }
    6  = @{ D = 44334; L = "Primera"; M = 100; N = 11000; O = 12000; P = 11500; Q = "`$/caja 12 kilos granel";     R = "Región de O'Higgins";  S = 11500; T = 1  }
    7  = @{ D = 44708; L = "Primera"; M = 70;  N = 12000; O = 13000; P = 12571; Q = "`$/caja 12 kilos empedrada";  R = "Provincia de Curicó";  S = 1048;  T = 12 }
    8  = @{ D = 44330; L = "Primera"; M = 100; N = 15000; O = 16000; P = 15500; Q = "`$/caja 18 kilos granel";     R = "Provincia de Curicó";  S = 861;   T = 18 }
    9  = @{ D = 44719; L = "Primera"; M = 50;  N = 14000; O = 15000; P = 14400; Q = "`$/caja 18 kilos granel";     R = "Región del Maule";     S = 800;   T = 18 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value  = $vals.D    # D: Fecha
    $ws.Cells.Item($row, 12).Value = $vals.L    # L: Calidad
    $ws.Cells.Item($row, 13).Value = $vals.M    # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals.N    # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals.O    # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals.P    # P: Precio promedio ponderado
    $ws.Cells.Item($row, 17).Value = $vals.Q    # Q: Unidad de comercializacion
    $ws.Cells.Item($row, 18).Value = $vals.R    # R: Origen
    $ws.Cells.Item($row, 19).Value = $vals.S    # S: Precio $/Kg
    $ws.Cells.Item($row, 20).Value = $vals.T    # T: Kg / unidad
}
